$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each D/E cell stores a literal numeric- or percent-looking string (the sheet
# uses text cells throughout, e.g. D2 = "276.13", E2 = "0.69%"). Setting
# .NumberFormat = "@" (Text) before writing .Value keeps Excel from silently
# auto-converting the new text into a Number/Percentage cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "275.99"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.63%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.13"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.63%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.857"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.05%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06400"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.16%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.943"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.80%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.198"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-5.41%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8783"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.99%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1519"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "4.17%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.05098"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.06%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07515"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.90%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.02975"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.74%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.08980"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.67%"

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-1.03%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006387"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.17%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006185"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "4.30%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.483"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.02%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.306"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.48%"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.86%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1348"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.924"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.70%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04417"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.65%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001176"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.02%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.003862"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-9.42%"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001200"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.07%"

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "14.73%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04161"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.93%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006811"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "2.18%"

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.46%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002150"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "2.94%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01186"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-2.81%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005276"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-0.16%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.680"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "15.75%"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-7.45%"
